$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# The "CasesTab" query (cell B2) had an erroneous extra `co:cohort` / `Cohort`
# return column appended to it (a leftover/mistaken variable from another
# query). Fix it by rewriting the query text without the Cohort clause and
# without the stray trailing blank line.
$casesQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n" +
              "WHERE demo.breed IN ['Weimaraner']`n" +
              "MATCH (c)<--(diag:diagnosis)`n" +
              "OPTIONAL MATCH (samp:sample)-->(c)`n" +
              "OPTIONAL MATCH (co:cohort)<-[*]-(c)`n" +
              "WITH DISTINCT c, s, demo, diag, co`n" +
              "RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n" +
              "        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n" +
              "        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n" +
              "        coalesce(demo.breed, '') AS Breed ,`n" +
              "        coalesce(diag.disease_term, '') AS Diagnosis ,`n" +
              "        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n" +
              "        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n" +
              "        coalesce(demo.sex, '') AS Sex ,`n" +
              "        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n" +
              "        coalesce(demo.weight, '') AS ``Weight (kg)``,`n" +
              "        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $casesQuery

# With the Cohort line removed the cell has fewer wrapped lines, so the
# (wrap-text) rows shrink to their new auto-fit heights.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 216
$ws.Rows.Item(4).RowHeight = 216

# Restore the view to the top of the sheet and select the query cell that
# was just fixed (previously the window was scrolled down with B4 selected).
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("B2").Select()
